$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2238.5
$ws.Cells.Item(19, 9).Value = 1595.5
$ws.Cells.Item(19, 11).Value = 1595.5
$ws.Cells.Item(19, 13).Value = -1420.5
$ws.Cells.Item(32, 8).Value = 3863.5
$ws.Cells.Item(32, 9).Value = 1260.3334
$ws.Cells.Item(32, 11).Value = 1260.3334
$ws.Cells.Item(32, 13).Value = -934.3334
$ws.Cells.Item(62, 8).Value = 2016.5
$ws.Cells.Item(62, 9).Value = 2016.5
$ws.Cells.Item(62, 11).Value = 2016.5
$ws.Cells.Item(62, 13).Value = -1392.5
$ws.Cells.Item(64, 8).Value = 5496.75
$ws.Cells.Item(64, 10).Value = 5496.75
$ws.Cells.Item(64, 12).Value = 5496.75
$ws.Cells.Item(64, 14).Value = -5992.75
$ws.Cells.Item(65, 8).Value = 2016.5
$ws.Cells.Item(65, 9).Value = 2016.5
$ws.Cells.Item(65, 11).Value = 10082.5
$ws.Cells.Item(65, 13).Value = -6962.5
$ws.Cells.Item(67, 8).Value = 5496.75
$ws.Cells.Item(67, 10).Value = 5496.75
$ws.Cells.Item(67, 12).Value = 5496.75
$ws.Cells.Item(67, 14).Value = -7212.75
$ws.Cells.Item(80, 8).Value = 2726.5557
$ws.Cells.Item(80, 10).Value = 2737.3635
$ws.Cells.Item(80, 12).Value = 8212.0905
$ws.Cells.Item(80, 14).Value = -10208.0905
$ws.Cells.Item(83, 8).Value = 2726.5557
$ws.Cells.Item(83, 10).Value = 2737.3635
$ws.Cells.Item(83, 12).Value = 24636.2715
$ws.Cells.Item(83, 14).Value = -34620.2715
$ws.Cells.Item(116, 8).Value = 5633.3335
$ws.Cells.Item(116, 10).Value = 8900
$ws.Cells.Item(116, 12).Value = 8900
$ws.Cells.Item(116, 14).Value = -15784
$ws.Cells.Item(127, 8).Value = 1851.3334
$ws.Cells.Item(127, 10).Value = 1611
$ws.Cells.Item(127, 12).Value = 4833
$ws.Cells.Item(127, 14).Value = -14753
$ws.Cells.Item(137, 8).Value = 2196.2307
$ws.Cells.Item(137, 9).Value = 1839
$ws.Cells.Item(137, 10).Value = 3000
$ws.Cells.Item(137, 11).Value = 5517
$ws.Cells.Item(137, 12).Value = 9000
$ws.Cells.Item(137, 13).Value = -2967
$ws.Cells.Item(137, 14).Value = -14100
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2408734.2
$ws.Cells.Item(32, 9).Value = 3046133
$ws.Cells.Item(32, 10).Value = 779826.3
$ws.Cells.Item(32, 11).Value = 3046133
$ws.Cells.Item(32, 12).Value = 779826.3
$ws.Cells.Item(32, 13).Value = -3045846
$ws.Cells.Item(32, 14).Value = -780400.3
$ws.Cells.Item(61, 8).Value = 2711
$ws.Cells.Item(61, 9).Value = 2711
$ws.Cells.Item(61, 11).Value = 2711
$ws.Cells.Item(61, 13).Value = -2499
$ws.Cells.Item(63, 8).Value = 5772.091
$ws.Cells.Item(63, 9).Value = 5710.1
$ws.Cells.Item(63, 10).Value = 6392
$ws.Cells.Item(63, 11).Value = 5710.1
$ws.Cells.Item(63, 12).Value = 6392
$ws.Cells.Item(63, 13).Value = -5024.1
$ws.Cells.Item(63, 14).Value = -7764
$ws.Cells.Item(66, 8).Value = 5772.091
$ws.Cells.Item(66, 9).Value = 5710.1
$ws.Cells.Item(66, 10).Value = 6392
$ws.Cells.Item(66, 11).Value = 28550.5
$ws.Cells.Item(66, 12).Value = 31960
$ws.Cells.Item(66, 13).Value = -25118.5
$ws.Cells.Item(66, 14).Value = -38824
$ws.Cells.Item(102, 8).Value = 1898.5
$ws.Cells.Item(102, 9).Value = 1898.5
$ws.Cells.Item(102, 11).Value = 1898.5
$ws.Cells.Item(102, 13).Value = -276.5
$ws.Cells.Item(110, 8).Value = 1683
$ws.Cells.Item(110, 9).Value = 1683
$ws.Cells.Item(110, 11).Value = 1683
$ws.Cells.Item(110, 13).Value = 362
$ws.Cells.Item(132, 8).Value = 2628.3684
$ws.Cells.Item(132, 9).Value = 2613.2778
$ws.Cells.Item(132, 11).Value = 7839.8334
$ws.Cells.Item(132, 13).Value = -5309.8334
$ws.Cells.Item(136, 8).Value = 2711
$ws.Cells.Item(136, 9).Value = 2711
$ws.Cells.Item(136, 11).Value = 8133
$ws.Cells.Item(136, 13).Value = -5583
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 399.75
$ws.Cells.Item(94, 9).Value = 350
$ws.Cells.Item(94, 10).Value = 449.5
$ws.Cells.Item(94, 11).Value = 350
$ws.Cells.Item(94, 12).Value = 449.5
$ws.Cells.Item(94, 13).Value = 101
$ws.Cells.Item(94, 14).Value = -1351.5
$ws.Cells.Item(99, 8).Value = 1999.8
$ws.Cells.Item(99, 9).Value = 1999.8
$ws.Cells.Item(99, 11).Value = 1999.8
$ws.Cells.Item(99, 13).Value = -501.8
$ws.Cells.Item(134, 8).Value = 2621
$ws.Cells.Item(134, 9).Value = 2621
$ws.Cells.Item(134, 11).Value = 7863
$ws.Cells.Item(134, 13).Value = -5328
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2416.818
$ws.Cells.Item(16, 9).Value = 2799
$ws.Cells.Item(16, 10).Value = 1748
$ws.Cells.Item(16, 11).Value = 2799
$ws.Cells.Item(16, 12).Value = 1748
$ws.Cells.Item(16, 13).Value = -2512
$ws.Cells.Item(16, 14).Value = -2322
$ws.Cells.Item(31, 8).Value = 1260.5
$ws.Cells.Item(31, 9).Value = 1213.25
$ws.Cells.Item(31, 11).Value = 1213.25
$ws.Cells.Item(31, 13).Value = -918.25
$ws.Cells.Item(34, 8).Value = 1260.5
$ws.Cells.Item(34, 9).Value = 1213.25
$ws.Cells.Item(34, 11).Value = 1213.25
$ws.Cells.Item(34, 13).Value = -1011.25
$ws.Cells.Item(58, 8).Value = 2035.5
$ws.Cells.Item(58, 9).Value = 1998.2
$ws.Cells.Item(58, 11).Value = 1998.2
$ws.Cells.Item(58, 13).Value = -1795.2
$ws.Cells.Item(107, 8).Value = 1159.909
$ws.Cells.Item(107, 9).Value = 1161.2858
$ws.Cells.Item(107, 10).Value = 1157.5
$ws.Cells.Item(107, 11).Value = 1161.2858
$ws.Cells.Item(107, 12).Value = 1157.5
$ws.Cells.Item(107, 13).Value = 758.7141999999999
$ws.Cells.Item(107, 14).Value = -4997.5
$ws.Cells.Item(113, 8).Value = 2416.818
$ws.Cells.Item(113, 9).Value = 2799
$ws.Cells.Item(113, 10).Value = 1748
$ws.Cells.Item(113, 11).Value = 2799
$ws.Cells.Item(113, 12).Value = 1748
$ws.Cells.Item(113, 13).Value = -629
$ws.Cells.Item(113, 14).Value = -6088
$ws.Cells.Item(134, 8).Value = 2124.6155
$ws.Cells.Item(134, 9).Value = 2120
$ws.Cells.Item(134, 11).Value = 6360
$ws.Cells.Item(134, 13).Value = -3825
$ws.Cells.Item(136, 8).Value = 2035.5
$ws.Cells.Item(136, 9).Value = 1998.2
$ws.Cells.Item(136, 11).Value = 5994.6
$ws.Cells.Item(136, 13).Value = -3444.6
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2530.923
$ws.Cells.Item(102, 9).Value = 2530.923
$ws.Cells.Item(102, 11).Value = 2530.923
$ws.Cells.Item(102, 13).Value = -908.9229999999998
$ws.Cells.Item(132, 8).Value = 2999.6667
$ws.Cells.Item(132, 9).Value = 2600
$ws.Cells.Item(132, 11).Value = 7800
$ws.Cells.Item(132, 13).Value = -5270
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(17, 8).Value = 122
$ws.Cells.Item(17, 9).Value = 8
$ws.Cells.Item(17, 10).Value = 350
$ws.Cells.Item(17, 11).Value = 8
$ws.Cells.Item(17, 12).Value = 350
$ws.Cells.Item(17, 13).Value = 162
$ws.Cells.Item(17, 14).Value = -690
$ws.Cells.Item(61, 8).Value = 2914.875
$ws.Cells.Item(61, 9).Value = 2902.7144
$ws.Cells.Item(61, 11).Value = 2902.7144
$ws.Cells.Item(61, 13).Value = -2700.7144
$ws.Cells.Item(93, 8).Value = 3321.8572
$ws.Cells.Item(93, 9).Value = 3000.6
$ws.Cells.Item(93, 11).Value = 3000.6
$ws.Cells.Item(93, 13).Value = -1752.6
$ws.Cells.Item(100, 8).Value = 3677.5
$ws.Cells.Item(100, 9).Value = 3206.6667
$ws.Cells.Item(100, 11).Value = 3206.6667
$ws.Cells.Item(100, 13).Value = -2665.6667
$ws.Cells.Item(105, 8).Value = 270307.5
$ws.Cells.Item(105, 10).Value = 270307.5
$ws.Cells.Item(105, 12).Value = 270307.5
$ws.Cells.Item(105, 14).Value = -277295.5
$ws.Cells.Item(113, 8).Value = 2914.875
$ws.Cells.Item(113, 9).Value = 2902.7144
$ws.Cells.Item(113, 11).Value = 2902.7144
$ws.Cells.Item(113, 13).Value = -732.7143999999998
$ws.Cells.Item(122, 8).Value = 6320.839
$ws.Cells.Item(122, 9).Value = 4927.857
$ws.Cells.Item(122, 10).Value = 7468
$ws.Cells.Item(122, 11).Value = 14783.571
$ws.Cells.Item(122, 12).Value = 22404
$ws.Cells.Item(122, 13).Value = -12333.571
$ws.Cells.Item(122, 14).Value = -27304
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 89547.5
$ws.Cells.Item(70, 9).Value = 90095
$ws.Cells.Item(70, 10).Value = 89000
$ws.Cells.Item(70, 11).Value = 90095
$ws.Cells.Item(70, 12).Value = 89000
$ws.Cells.Item(70, 13).Value = -89780
$ws.Cells.Item(70, 14).Value = -89630
$ws.Cells.Item(73, 8).Value = 89547.5
$ws.Cells.Item(73, 9).Value = 90095
$ws.Cells.Item(73, 10).Value = 89000
$ws.Cells.Item(73, 11).Value = 90095
$ws.Cells.Item(73, 12).Value = 89000
$ws.Cells.Item(73, 13).Value = -89003
$ws.Cells.Item(73, 14).Value = -91184
$ws.Cells.Item(105, 8).Value = 7450
$ws.Cells.Item(105, 10).Value = 7450
$ws.Cells.Item(105, 12).Value = 7450
$ws.Cells.Item(105, 14).Value = -14438
$ws.Cells.Item(126, 8).Value = 6244.9653
$ws.Cells.Item(126, 9).Value = 5209.25
$ws.Cells.Item(126, 10).Value = 7519.6924
$ws.Cells.Item(126, 11).Value = 15627.75
$ws.Cells.Item(126, 12).Value = 22559.0772
$ws.Cells.Item(126, 13).Value = -13157.75
$ws.Cells.Item(126, 14).Value = -27499.0772
